$wb = $excel.ActiveWorkbook

# --- Sheet: Weekly Quantity --- append rows 71-72
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Cells.Item(71, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(71, 1).Value = 45669.99999999999
$ws1.Cells.Item(71, 2).Value = 14
$ws1.Cells.Item(72, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(72, 1).Value = 45676.99999999999
$ws1.Cells.Item(72, 2).Value = 33

# --- Sheet: Monthly Trend --- append row 26
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Cells.Item(26, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(26, 1).Value = 45688.99999999999
$ws2.Cells.Item(26, 2).Value = 47

# --- Sheet: PO Forecast --- new forecast model for rows 2-80 ---
$ws3 = $wb.Worksheets.Item("PO Forecast")
$ws3.Cells.Item(2, 1).Value = 44934.99999999999
$ws3.Cells.Item(2, 2).Value = 301
$ws3.Cells.Item(3, 1).Value = 44948.99999999999
$ws3.Cells.Item(3, 2).Value = 66
$ws3.Cells.Item(4, 1).Value = 44969.99999999999
$ws3.Cells.Item(4, 2).Value = 156
$ws3.Cells.Item(5, 1).Value = 44976.99999999999
$ws3.Cells.Item(5, 2).Value = 168
$ws3.Cells.Item(6, 1).Value = 44983.99999999999
$ws3.Cells.Item(6, 2).Value = 146
$ws3.Cells.Item(7, 1).Value = 44990.99999999999
$ws3.Cells.Item(7, 2).Value = 124
$ws3.Cells.Item(8, 1).Value = 45011.99999999999
$ws3.Cells.Item(8, 2).Value = 136
$ws3.Cells.Item(9, 1).Value = 45018.99999999999
$ws3.Cells.Item(9, 2).Value = 139
$ws3.Cells.Item(10, 1).Value = 45025.99999999999
$ws3.Cells.Item(10, 2).Value = 150
$ws3.Cells.Item(11, 1).Value = 45032.99999999999
$ws3.Cells.Item(11, 2).Value = 169
$ws3.Cells.Item(12, 1).Value = 45053.99999999999
$ws3.Cells.Item(12, 2).Value = 84
$ws3.Cells.Item(13, 1).Value = 45060.99999999999
$ws3.Cells.Item(13, 2).Value = 45
$ws3.Cells.Item(14, 1).Value = 45067.99999999999
$ws3.Cells.Item(14, 2).Value = 53
$ws3.Cells.Item(15, 1).Value = 45074.99999999999
$ws3.Cells.Item(15, 2).Value = 87
$ws3.Cells.Item(16, 1).Value = 45081.99999999999
$ws3.Cells.Item(16, 2).Value = 97
$ws3.Cells.Item(17, 1).Value = 45088.99999999999
$ws3.Cells.Item(17, 2).Value = 63
$ws3.Cells.Item(18, 1).Value = 45095.99999999999
$ws3.Cells.Item(18, 2).Value = 26
$ws3.Cells.Item(19, 1).Value = 45102.99999999999
$ws3.Cells.Item(19, 2).Value = 46
$ws3.Cells.Item(20, 1).Value = 45109.99999999999
$ws3.Cells.Item(20, 2).Value = 128
$ws3.Cells.Item(21, 1).Value = 45116.99999999999
$ws3.Cells.Item(21, 2).Value = 202
$ws3.Cells.Item(22, 1).Value = 45123.99999999999
$ws3.Cells.Item(22, 2).Value = 193
$ws3.Cells.Item(23, 1).Value = 45130.99999999999
$ws3.Cells.Item(23, 2).Value = 104
$ws3.Cells.Item(24, 1).Value = 45137.99999999999
$ws3.Cells.Item(24, 2).Value = 29
$ws3.Cells.Item(25, 1).Value = 45144.99999999999
$ws3.Cells.Item(25, 2).Value = 61
$ws3.Cells.Item(26, 1).Value = 45151.99999999999
$ws3.Cells.Item(26, 2).Value = 191
$ws3.Cells.Item(27, 1).Value = 45158.99999999999
$ws3.Cells.Item(27, 2).Value = 303
$ws3.Cells.Item(28, 1).Value = 45165.99999999999
$ws3.Cells.Item(28, 2).Value = 286
$ws3.Cells.Item(29, 1).Value = 45172.99999999999
$ws3.Cells.Item(29, 2).Value = 150
$ws3.Cells.Item(30, 1).Value = 45179.99999999999
$ws3.Cells.Item(30, 2).Value = 21
$ws3.Cells.Item(31, 1).Value = 45186.99999999999
$ws3.Cells.Item(31, 2).Value = 21
$ws3.Cells.Item(32, 1).Value = 45193.99999999999
$ws3.Cells.Item(32, 2).Value = 150
$ws3.Cells.Item(33, 1).Value = 45207.99999999999
$ws3.Cells.Item(33, 2).Value = 304
$ws3.Cells.Item(34, 1).Value = 45214.99999999999
$ws3.Cells.Item(34, 2).Value = 190
$ws3.Cells.Item(35, 1).Value = 45221.99999999999
$ws3.Cells.Item(35, 2).Value = 45
$ws3.Cells.Item(36, 1).Value = 45228.99999999999
$ws3.Cells.Item(36, 2).Value = 0
$ws3.Cells.Item(37, 1).Value = 45235.99999999999
$ws3.Cells.Item(37, 2).Value = 34
$ws3.Cells.Item(38, 1).Value = 45242.99999999999
$ws3.Cells.Item(38, 2).Value = 115
$ws3.Cells.Item(39, 1).Value = 45249.99999999999
$ws3.Cells.Item(39, 2).Value = 138
$ws3.Cells.Item(40, 1).Value = 45256.99999999999
$ws3.Cells.Item(40, 2).Value = 91
$ws3.Cells.Item(41, 1).Value = 45263.99999999999
$ws3.Cells.Item(41, 2).Value = 43
$ws3.Cells.Item(42, 1).Value = 45270.99999999999
$ws3.Cells.Item(42, 2).Value = 76
$ws3.Cells.Item(43, 1).Value = 45277.99999999999
$ws3.Cells.Item(43, 2).Value = 197
$ws3.Cells.Item(44, 1).Value = 45298.99999999999
$ws3.Cells.Item(44, 2).Value = 318
$ws3.Cells.Item(45, 1).Value = 45305.99999999999
$ws3.Cells.Item(45, 2).Value = 185
$ws3.Cells.Item(46, 1).Value = 45312.99999999999
$ws3.Cells.Item(46, 2).Value = 74
$ws3.Cells.Item(47, 1).Value = 45319.99999999999
$ws3.Cells.Item(47, 2).Value = 46
$ws3.Cells.Item(48, 1).Value = 45326.99999999999
$ws3.Cells.Item(48, 2).Value = 89
$ws3.Cells.Item(49, 1).Value = 45333.99999999999
$ws3.Cells.Item(49, 2).Value = 145
$ws3.Cells.Item(50, 1).Value = 45340.99999999999
$ws3.Cells.Item(50, 2).Value = 165
$ws3.Cells.Item(51, 1).Value = 45347.99999999999
$ws3.Cells.Item(51, 2).Value = 147
$ws3.Cells.Item(52, 1).Value = 45354.99999999999
$ws3.Cells.Item(52, 2).Value = 123
$ws3.Cells.Item(53, 1).Value = 45361.99999999999
$ws3.Cells.Item(53, 2).Value = 116
$ws3.Cells.Item(54, 1).Value = 45368.99999999999
$ws3.Cells.Item(54, 2).Value = 125
$ws3.Cells.Item(55, 1).Value = 45375.99999999999
$ws3.Cells.Item(55, 2).Value = 132
$ws3.Cells.Item(56, 1).Value = 45382.99999999999
$ws3.Cells.Item(56, 2).Value = 134
$ws3.Cells.Item(57, 1).Value = 45389.99999999999
$ws3.Cells.Item(57, 2).Value = 143
$ws3.Cells.Item(58, 1).Value = 45417.99999999999
$ws3.Cells.Item(58, 2).Value = 90
$ws3.Cells.Item(59, 1).Value = 45431.99999999999
$ws3.Cells.Item(59, 2).Value = 45
$ws3.Cells.Item(60, 1).Value = 45445.99999999999
$ws3.Cells.Item(60, 2).Value = 95
$ws3.Cells.Item(61, 1).Value = 45452.99999999999
$ws3.Cells.Item(61, 2).Value = 68
$ws3.Cells.Item(62, 1).Value = 45487.99999999999
$ws3.Cells.Item(62, 2).Value = 199
$ws3.Cells.Item(63, 1).Value = 45494.99999999999
$ws3.Cells.Item(63, 2).Value = 119
$ws3.Cells.Item(64, 1).Value = 45536.99999999999
$ws3.Cells.Item(64, 2).Value = 175
$ws3.Cells.Item(65, 1).Value = 45543.99999999999
$ws3.Cells.Item(65, 2).Value = 33
$ws3.Cells.Item(66, 1).Value = 45578.99999999999
$ws3.Cells.Item(66, 2).Value = 213
$ws3.Cells.Item(67, 1).Value = 45592.99999999999
$ws3.Cells.Item(67, 2).Value = 0
$ws3.Cells.Item(68, 1).Value = 45613.99999999999
$ws3.Cells.Item(68, 2).Value = 137
$ws3.Cells.Item(69, 1).Value = 45620.99999999999
$ws3.Cells.Item(69, 2).Value = 98
$ws3.Cells.Item(70, 1).Value = 45641.99999999999
$ws3.Cells.Item(70, 2).Value = 167
$ws3.Cells.Item(71, 1).Value = 45669.99999999999
$ws3.Cells.Item(71, 2).Value = 206
$ws3.Cells.Item(72, 1).Value = 45676.99999999999
$ws3.Cells.Item(72, 2).Value = 85
$ws3.Cells.Item(73, 1).Value = 45683.99999999999
$ws3.Cells.Item(73, 2).Value = 41
$ws3.Cells.Item(74, 1).Value = 45690.99999999999
$ws3.Cells.Item(74, 2).Value = 75
$ws3.Cells.Item(75, 1).Value = 45697.99999999999
$ws3.Cells.Item(75, 2).Value = 133
$ws3.Cells.Item(76, 1).Value = 45704.99999999999
$ws3.Cells.Item(76, 2).Value = 161
$ws3.Cells.Item(77, 1).Value = 45711.99999999999
$ws3.Cells.Item(77, 2).Value = 149
$ws3.Cells.Item(78, 1).Value = 45718.99999999999
$ws3.Cells.Item(78, 2).Value = 123
$ws3.Cells.Item(79, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(79, 1).Value = 45725.99999999999
$ws3.Cells.Item(79, 2).Value = 112
$ws3.Cells.Item(80, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(80, 1).Value = 45732.99999999999
$ws3.Cells.Item(80, 2).Value = 119
